$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 88: the date/time serial in A88 was recomputed by the R script ---
# (keeps its existing date-time style; only the value changes)
$ws.Range("A88").Value = 45463.2916666667

# --- Row 89: a new OHLCV data point appended by the R script ---
# Copy row 88's formatting down first (gives A89 the date style and H89 the
# "PAL.MI" ticker text already), then overwrite every cell with the new data.
$ws.Range("A88:H88").Copy($ws.Range("A89:H89"))

$ws.Range("A89").Value = 45464.6140509259
$ws.Range("B89").Value = 6900
$ws.Range("C89").Value = 6.19999980926514
$ws.Range("D89").Value = 6
$ws.Range("E89").Value = 6.07999992370605
$ws.Range("F89").Value = 6.19999980926514

# G89 ("adj_close") is stored as text equal to the close price, "6.19999980926514".
# A bare assignment would be auto-parsed back into a number, so force text entry
# the same way Excel's UI does (leading apostrophe), then drop the resulting
# "quote prefix" direct formatting so the cell is plain text with no style override.
$ws.Range("G89").Value = "'6.19999980926514"
$ws.Range("G89").Style = $ws.Range("B89").Style

# H89 ("ticker") already reads "PAL.MI" from the row-88 copy, matching the source data.
